# Generate Report for Handoff
#
# The "a694d871-5840-414a-83f6-5f17fd7be361" source file has dropped out of
# the report (handed off + closed out), while the "8bc2bcd8-..." file has
# moved from "Handed back: in sync with en-US" to "Ready for handoff" with a
# refreshed handoff timestamp. Apply this on all three worksheets:
#   - Overview : drop the a694d871 row, update the status cells
#   - zh-cn    : drop the a694d871 row, update status + handoff datetime
#   - de-de    : drop the a694d871 row, update status + handoff datetime

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview  (columns: File Name | zh-cn | de-de)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"

# Snapshot the hyperlinks (in sheet order) before we start mutating rows, so
# we keep stable references regardless of how the live collection re-indexes.
$links = @()
foreach ($hl in $ws.Hyperlinks) { $links += $hl }
# links: 0 -> A2 (8bc2bcd8.md), 1 -> A3 (a694d871.md), 2 -> A4 (.localization-config)

# The row-3 hyperlink is about to become the (shifted-up) .localization-config
# row, so repoint it before deleting the now-redundant row-4 hyperlink.
$links[1].Address = $links[2].Address
$links[1].TextToDisplay = $links[2].TextToDisplay
$links[2].Delete()

$ws.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-09 05:49:12"

$links = @()
foreach ($hl in $ws.Hyperlinks) { $links += $hl }
# links: 0:A2 1:C2 2:E2 3:F2 (8bc2bcd8 row) | 4:A3 5:C3 6:E3 7:F3 (a694d871 row) | 8:A4 (.localization-config)

$links[4].Address = $links[8].Address
$links[4].TextToDisplay = $links[8].TextToDisplay

# Delete the stale links highest-index-first so earlier snapshot refs stay valid.
$links[8].Delete()
$links[7].Delete()
$links[6].Delete()
$links[5].Delete()

$ws.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-09 05:49:20"

$links = @()
foreach ($hl in $ws.Hyperlinks) { $links += $hl }
# links: 0:A2 1:C2 2:E2 3:F2 (8bc2bcd8 row) | 4:A3 5:C3 6:E3 7:F3 (a694d871 row) | 8:A4 (.localization-config)

$links[4].Address = $links[8].Address
$links[4].TextToDisplay = $links[8].TextToDisplay

$links[8].Delete()
$links[7].Delete()
$links[6].Delete()
$links[5].Delete()

$ws.Rows.Item(3).Delete()
